$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.893.51"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "1.552.10"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'206.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'21.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").Value = "'0.0586"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "'0.0859"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "1.772.58"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").Value = "1.553.88"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "26.877.62"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'61.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'216.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.19%  "
$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "'7.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("D25").Value = "'153.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").Value = "'14.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Value = "1.418.14"
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("E35").Value = "  +3.77%  "
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").Value = "'0.0165"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "'0.807"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").Value = "'0.988"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("D45").Value = "'63.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").Value = "'1.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "1.686.75"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "'86.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "'0.0518"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").Value = "0.0₇0975"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "'0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.70%  "
